# Bugfixes and general improvements
# The "livestock_area" grouping-variable row (varName=livestock_area,
# label=Land Area for Livestock, shortName=Farm Size, Levels=0,1,2,3,
# Labels="0 ha,>0-2 ha,>2-4 ha,>4 ha", level=All) is removed from the
# lookup table. Deleting the entire worksheet row shifts every row below
# it up by one and drops the now-unused shared strings automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 holds the "livestock_area" / "Farm Size" entry - remove it entirely.
$ws.Rows("5").Delete()

# Match the author's final selection/active cell.
$ws.Range("C16").Select()
